# Update the "Pais" worksheet with the refreshed COVID figures (3 Jul 2020, 19:49)
# and the handful of countries whose ranking in the total-cases sort swapped
# (Suazilandia/Libia, Comoras/Birmania, Dominica/Fiyi, Islas Malvinas/Groenlandia).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 1;   A = 'Datos actualizados a 3 de Julio de 2020 a las 19:49' },
    @{ Row = 4;   B = 2860886; C = 25202; D = 1199235; E = 1529844; G = 322; H = 131807 },
    @{ Row = 5;   B = 1508991; C = 7638;  E = 530540;  G = 314; H = 62304 },
    @{ Row = 7;   B = 647503;  C = 20335; D = 392882;  E = 235960; G = 436; H = 18661 },
    @{ Row = 16;  B = 203456;  C = 1172;  D = 178278;  E = 19992;  G = 19;  H = 5186 },
    @{ Row = 18;  B = 196823;  C = 106;   E = 6759 },
    @{ Row = 23;  B = 105025;  C = 253;   D = 68650;   E = 27712;  G = 21;  H = 8663 },
    @{ Row = 46;  B = 35148;   C = 951;   D = 18392;   E = 15981;  G = 10;  H = 775 },
    @{ Row = 54;  B = 25498;   C = 9;     E = 394;     G = 2;      H = 1740 },
    @{ Row = 61;  B = 17445;   C = 295;   E = 6780;    G = 12;     H = 572 },
    @{ Row = 64;  B = 15070;   C = 413;   D = 10832;   E = 3301;   G = 9;   H = 937 },
    @{ Row = 65;  B = 13288;   C = 319;   D = 9160;    E = 3898 },
    @{ Row = 73;  B = 8916;    C = 14;    E = 527 },
    @{ Row = 94;  B = 4447;    C = 52;    D = 4016;    E = 321 },
    @{ Row = 100; E = 2674;    G = 2;     H = 11 },
    @{ Row = 108; B = 2410;    C = 10;    D = 1976;    E = 424 },
    @{ Row = 109; B = 2361;    C = 8;     D = 2224;    E = 51 },
    @{ Row = 111; B = 2285;    C = 25;    D = 1507;    E = 661 },
    @{ Row = 123; B = 1524;    C = 6;     D = 1042;    E = 420;    G = 2;   H = 62 },
    @{ Row = 128; B = 1240;    C = 19;    D = 536;     E = 369;    G = 10;  H = 335 },
    @{ Row = 131; E = 240;     G = 1;     H = 10 },
    @{ Row = 141; A = 'Suazilandia'; B = 909; C = 36; D = 515; E = 381; G = 2; H = 13 },
    @{ Row = 142; A = 'Libia';       B = 891;          D = 224; E = 641;          H = 26 },
    @{ Row = 143; B = 871;     C = 3;     D = 786;     E = 11 },
    @{ Row = 145; B = 833;     C = 14;    D = 346;     E = 450 },
    @{ Row = 164; A = 'Comoras';  B = 309; C = 6; D = 241; E = 61; H = 7 },
    @{ Row = 165; A = 'Birmania'; B = 306; C = 2; D = 237; E = 63; H = 6 },
    @{ Row = 205; A = 'Dominica' },
    @{ Row = 206; A = 'Fiyi' },
    @{ Row = 209; A = 'Islas Malvinas' },
    @{ Row = 210; A = 'Groenlandia' }
)

foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @('A', 'B', 'C', 'D', 'E', 'F', 'G', 'H')) {
        if ($u.ContainsKey($col)) {
            $ws.Range("$col$r").Value = $u[$col]
        }
    }
}
